$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Every slide's "date" placeholder was re-stamped from 2020/9/24
#    to 2020/9/25 (cosmetic footer date shown on every slide).
# ------------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $placeholders = $slide.Shapes.Placeholders
    for ($j = 1; $j -le $placeholders.Count; $j++) {
        $shape = $placeholders.Item($j)
        if ($shape.PlaceholderFormat.Type -eq 16) {
            $dateRange = $shape.TextFrame.TextRange
            if ($dateRange.Text -eq "2020/9/24") {
                $dateRange.Text = "2020/9/25"
            }
        }
    }
}

# ------------------------------------------------------------------
# 2) Slide 2 ("5 Install TensorFlow") content tweaks: the YouTube
#    call-out line now says it is specifically for MacOS, and the
#    link text underneath points at the new, shorter mobile URL.
# ------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$subtitle = $slide2.Shapes.Item(2)
$bodyRange = $subtitle.TextFrame.TextRange

$label = $bodyRange.Find("YouTube video:")
if ($label -ne $null) {
    $label.Text = "YouTube video for MacOS:"
}

$link = $bodyRange.Find("https://www.youtube.com/watch?v=RgO8BBNGB8w&t=376s")
if ($link -ne $null) {
    $link.Text = "https://m.youtube.com/watch?v=MpUvdLD932c"
}
